$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet named "2022-Q3" right after "总计" and
#    before "2022-Q2" (i.e. it becomes the new 2nd sheet).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

# Source sheet to copy cell formatting from (the sheet that is currently
# named "2022-Q2", used as a template for header / index-column styling).
$template = $wb.Worksheets.Item("2022-Q2")

# Copy the header-row formatting (bold, bordered, centered style) from the
# template sheet onto the new sheet's header row.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy the index-column (column A) formatting used for the numeric row index.
$template.Range("A2").Copy()
$newSheet.Range("A2:A19").PasteSpecial(-4122)

# Header labels.
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Fund holding rows: code, name, scale, stock position, position pct,
# market value (亿元), position rank.
$rows = @(
    @("100061", "富国中国中小盘混合（QDII）人民币", "35.11", "83.32", "6.50", "2.2822", 2),
    @("010591", "富国中国中小盘混合（QDII）美元",   "35.11", "83.32", "6.50", "2.2822", 2),
    @("011006", "工银圆丰三年持有期混合",           "68.77", "88.83", "2.68", "1.8430", 8),
    @("009076", "工银圆兴混合",                     "49.75", "87.73", "3.68", "1.8308", 6),
    @("007139", "富国民裕进取沪港深成长精选混合A",   "10.56", "88.09", "5.07", "0.5354", 9),
    @("009029", "工银高质量成长混合A",               "12.13", "86.00", "3.58", "0.4343", 6),
    @("006752", "天弘港股通精选灵活配置混合A",       "4.67",  "93.60", "7.93", "0.3703", 2),
    @("006753", "天弘港股通精选灵活配置混合C",       "2.52",  "93.60", "7.93", "0.1998", 2),
    @("011556", "富国民裕进取沪港深成长精选混合C",   "2.46",  "88.09", "5.07", "0.1247", 9),
    @("012584", "南方中国新兴经济9个月持有期混合（QDII）A", "2.69", "91.51", "3.49", "0.0939", 8),
    @("009240", "泰康蓝筹优势一年持有期股票",         "4.69",  "71.86", "1.56", "0.0732", 10),
    @("007109", "南方沪港深核心优势混合",             "1.59",  "85.13", "4.04", "0.0642", 7),
    @("009030", "工银高质量成长混合C",                 "1.64",  "86.00", "3.58", "0.0587", 6),
    @("005259", "建信龙头企业股票",                   "0.76",  "84.43", "4.81", "0.0366", 6),
    @("011969", "建信港股通精选混合A",                 "0.56",  "63.37", "5.02", "0.0281", 7),
    @("011970", "建信港股通精选混合C",                 "0.23",  "63.37", "5.02", "0.0115", 7),
    @("004098", "前海开源港股通股息率50强股票",       "0.31",  "87.14", "2.40", "0.0074", 10),
    @("012585", "南方中国新兴经济9个月持有期混合（QDII）C", "0.10", "91.51", "3.49", "0.0035", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row 2 holding the
#    2022-Q3 totals, pushing all existing quarters down by one row.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 18
$total.Cells.Item(2, 4).Value = 10.28

# The row-insert operation copies down the formatting of the row above
# (the header), so strip that back out of the plain data cells B2:D2.
$total.Range("B2:D2").Style = "Normal"

# Re-apply the index-column style (also lost on insert) by copying it from
# the row below, which already carries the correct style.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$total.Cells.Item(2, 1).Value = 0
